$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.491.86"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").Value = "1.977.35"
$ws.Range("E3").Value = "  -3.60%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'245.36"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  -4.60%  "

$ws.Range("E7").Value = "  +5.09%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'58.44"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "'0.358"
$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("D13").Value = "'0.942"
$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("D14").Value = "'14.39"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").Value = "2.262.04"
$ws.Range("E15").Value = "  -3.83%  "

$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("D17").Value = "1.968.98"
$ws.Range("E17").Value = "  -3.91%  "

$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "  +6.12%  "

$ws.Range("D19").Value = "35.504.07"
$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("D20").Value = "'71.26"
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("D22").Value = "'232.26"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  +20.43%  "

$ws.Range("E26").Value = "  -2.77%  "

$ws.Range("D27").Value = "'164.19"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  -2.99%  "

$ws.Range("D29").Value = "'19.08"
$ws.Range("E29").Value = "  -4.60%  "

$ws.Range("E30").Value = "  -2.74%  "

$ws.Range("E31").Value = "  -3.65%  "

$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  -8.37%  "

$ws.Range("D33").Value = "'0.0940"
$ws.Range("E33").Value = "  +15.33%  "

$ws.Range("D34").Value = "'0.0590"
$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("E36").Value = "  +8.84%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  -4.25%  "

$ws.Range("E39").Value = "  +7.40%  "

$ws.Range("E40").Value = "  -1.85%  "

$ws.Range("D41").Value = "'2.88"
$ws.Range("E41").Value = "  +1.91%  "

$ws.Range("E42").Value = "  -2.17%  "

$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("D44").Value = "'91.13"
$ws.Range("E44").Value = "  -2.44%  "

$ws.Range("D45").Value = "'15.91"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").Value = "'0.0882"
$ws.Range("E46").Value = "  -5.68%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.50"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.365.23"
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("E49").Value = "  +1.29%  "

$ws.Range("D50").Value = "'46.71"
$ws.Range("E50").Value = "  +4.21%  "

$ws.Range("D51").Value = "'3.67"
$ws.Range("E51").Value = "  +14.12%  "

